$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "C13" component row (row 14) -- per commit message, this
# designator is not necessary (duplicate near CP2102). Deleting the
# entire row shifts every row below it up by one.
$ws.Rows.Item(14).Delete()

# Mirror Excel's behavior of leaving the row that now occupies the
# deleted row's position selected (whole row 14 selected).
[void]$ws.Rows.Item(14).Select()
